$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Save" in column H, row 1, matching format/style of existing header (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for rows 2-5 (plain, unstyled numeric cells)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
